$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 299, shifting the existing data
# (old rows 299..406) down to rows 301..408.
$ws.Rows("299:300").Insert()

# Populate the newly-inserted row 299 with a new "Coliflor" record
# (Terminal Hortofruticola Agro Chillan, calidad Primera).
$ws.Range("A299").Value = 7
$ws.Range("B299").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C299").Value = "Ñuble"
$ws.Range("D299").Value = 44988
$ws.Range("E299").Value = 16
$ws.Range("F299").Value = 100112008
$ws.Range("G299").Value = "Coliflor"
$ws.Range("H299").Value = "Sin especificar"
$ws.Range("I299").Value = "Primera"
$ws.Range("J299").Value = 300
$ws.Range("K299").Value = 1000
$ws.Range("L299").Value = 1000
$ws.Range("M299").Value = 1000
$ws.Range("N299").Value = "$/unidad"
$ws.Range("O299").Value = "Región del Maule"
$ws.Range("P299").Value = 1000
$ws.Range("Q299").Value = 1
$ws.Range("R299").Value = "Hortaliza"

# Populate the newly-inserted row 300 with a new "Coliflor" record
# (Terminal Hortofruticola Agro Chillan, calidad Segunda).
$ws.Range("A300").Value = 7
$ws.Range("B300").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C300").Value = "Ñuble"
$ws.Range("D300").Value = 44988
$ws.Range("E300").Value = 16
$ws.Range("F300").Value = 100112008
$ws.Range("G300").Value = "Coliflor"
$ws.Range("H300").Value = "Sin especificar"
$ws.Range("I300").Value = "Segunda"
$ws.Range("J300").Value = 300
$ws.Range("K300").Value = 800
$ws.Range("L300").Value = 800
$ws.Range("M300").Value = 800
$ws.Range("N300").Value = "$/unidad"
$ws.Range("O300").Value = "Región del Maule"
$ws.Range("P300").Value = 800
$ws.Range("Q300").Value = 1
$ws.Range("R300").Value = "Hortaliza"
